$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7: SKU 0006, LED Bicyle Front Head Light, Landing 450, Selling 799, Units 1
$ws.Range("A7").Value = "0006"
$ws.Range("B7").Value = "LED Bicyle Front Head Light "
$ws.Range("C7").Value = 450
$ws.Range("D7").Value = 799
$ws.Range("E7").Value = 1

# Move the active selection to the newly added price cell, matching the author's edit
$ws.Range("C7").Select()
